$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 95, shifting the existing rows 95-127 down to 96-128.
$ws.Rows.Item(95).Insert()

# Populate the newly inserted row 95 with the new weekly price observation.
$ws.Range("A95").Value = 4
$ws.Range("B95").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C95").Value = "Los Lagos"
$ws.Range("D95").Value = 44985
$ws.Range("E95").Value = 10
$ws.Range("F95").Value = 100112031
$ws.Range("G95").Value = "Poroto verde"
$ws.Range("H95").Value = "Magnum"
$ws.Range("I95").Value = "Primera"
$ws.Range("J95").Value = 45
$ws.Range("K95").Value = 30000
$ws.Range("L95").Value = 30000
$ws.Range("M95").Value = 30000
$ws.Range("N95").Value = "`$/saco 25 kilos"
$ws.Range("O95").Value = "Región Metropolitana"
$ws.Range("P95").Value = 1200
$ws.Range("Q95").Value = 25
$ws.Range("R95").Value = "Hortaliza"
